$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snippets")

# Step 1: add Workbook.getActiveCell row (columns A-C first, D filled in later)
$ws.Cells.Item(110, 1).Value = "Workbook"
$ws.Cells.Item(110, 2).Value = "getActiveCell"
$ws.Cells.Item(110, 3).Value = "excel-workbook-get-active-cell"

# Step 2: add Worksheet.copy row (columns A-C first, D filled in later)
$ws.Cells.Item(111, 1).Value = "Worksheet"
$ws.Cells.Item(111, 2).Value = "copy"
$ws.Cells.Item(111, 3).Value = "excel-worksheet-copy"

# Step 3: Worksheet.protect (data in worksheet), fully filled in
$ws.Cells.Item(112, 1).Value = "Worksheet"
$ws.Cells.Item(112, 2).Value = "protect"
$ws.Cells.Item(112, 3).Value = "excel-protect-data-in-worksheet-and-workbook-structure"
$ws.Cells.Item(112, 4).Value = "protectDataInWorksheet"

# Step 4: Workbook.protect (workbook structure), fully filled in
$ws.Cells.Item(114, 1).Value = "Workbook"
$ws.Cells.Item(114, 2).Value = "protect"
$ws.Cells.Item(114, 3).Value = "excel-protect-data-in-worksheet-and-workbook-structure"
$ws.Cells.Item(114, 4).Value = "protectWorkbookStructure"

# Step 5: go back and fill in the deferred "run" snippet name for the first two rows
$ws.Cells.Item(110, 4).Value = "run"
$ws.Cells.Item(111, 4).Value = "run"

# Step 6: Worksheet.unprotect (data in worksheet)
$ws.Cells.Item(113, 1).Value = "Worksheet"
$ws.Cells.Item(113, 2).Value = "unprotect"
$ws.Cells.Item(113, 3).Value = "excel-protect-data-in-worksheet-and-workbook-structure"
$ws.Cells.Item(113, 4).Value = "unprotectDataInWorksheet"

# Step 7: Workbook.unprotect (workbook structure)
$ws.Cells.Item(115, 1).Value = "Workbook"
$ws.Cells.Item(115, 2).Value = "unprotect"
$ws.Cells.Item(115, 3).Value = "excel-protect-data-in-worksheet-and-workbook-structure"
$ws.Cells.Item(115, 4).Value = "unprotectWorkbookStructure"

# Step 8: Worksheet.protect (password-protected data in worksheet)
$ws.Cells.Item(116, 1).Value = "Worksheet"
$ws.Cells.Item(116, 2).Value = "protect"
$ws.Cells.Item(116, 3).Value = "excel-protect-data-in-worksheet-and-workbook-structure"
$ws.Cells.Item(116, 4).Value = "passwordProtectDataInWorksheet"

# Step 9: Worksheet.unprotect (password-protected data in worksheet)
$ws.Cells.Item(117, 1).Value = "Worksheet"
$ws.Cells.Item(117, 2).Value = "unprotect"
$ws.Cells.Item(117, 3).Value = "excel-protect-data-in-worksheet-and-workbook-structure"
$ws.Cells.Item(117, 4).Value = "passwordUnprotectDataInWorksheet"

# Step 10: Workbook.protect (password-protected workbook structure)
$ws.Cells.Item(118, 1).Value = "Workbook"
$ws.Cells.Item(118, 2).Value = "protect"
$ws.Cells.Item(118, 3).Value = "excel-protect-data-in-worksheet-and-workbook-structure"
$ws.Cells.Item(118, 4).Value = "passwordProtectWorkbookStructure"

# Step 11: Workbook.unprotect (password-protected workbook structure)
$ws.Cells.Item(119, 1).Value = "Workbook"
$ws.Cells.Item(119, 2).Value = "unprotect"
$ws.Cells.Item(119, 3).Value = "excel-protect-data-in-worksheet-and-workbook-structure"
$ws.Cells.Item(119, 4).Value = "passwordUnprotectWorkbookStructure"

# Resize the Snippets table to include the new rows
$table = $ws.ListObjects.Item("Snippets")
$table.Resize($ws.Range("A1:D119"))

# Update selection / frozen pane view to reflect the end-state
$ws.Range("A120").Select()
